$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.01949460275857007
$ws.Cells.Item(2, 4).Value = 0.1550633964863124
$ws.Cells.Item(2, 5).Value = 0.1571220482723845
$ws.Cells.Item(2, 6).Value = 1.708931743415746
$ws.Cells.Item(2, 7).Value = 1.075093835350231
$ws.Cells.Item(2, 8).Value = 1.050068151600698
$ws.Cells.Item(2, 9).Value = 0.918156576558296
$ws.Cells.Item(2, 10).Value = 0.1948255997740915
$ws.Cells.Item(2, 11).Value = 1.655386014733381
$ws.Cells.Item(2, 14).Value = 1.239475840486527

$ws.Cells.Item(3, 3).Value = 0.01832807991338825
$ws.Cells.Item(3, 4).Value = 0.1514779638615806
$ws.Cells.Item(3, 5).Value = 0.153602383405218
$ws.Cells.Item(3, 6).Value = 1.700992747160967
$ws.Cells.Item(3, 7).Value = 1.067889848414708
$ws.Cells.Item(3, 8).Value = 1.05283553585322
$ws.Cells.Item(3, 9).Value = 0.9129079596918359
$ws.Cells.Item(3, 10).Value = 0.1905101747098357
$ws.Cells.Item(3, 11).Value = 1.50418787687812
$ws.Cells.Item(3, 14).Value = 1.251296326750733

$ws.Cells.Item(4, 3).Value = 0.01760542479859595
$ws.Cells.Item(4, 4).Value = 0.149327481499455
$ws.Cells.Item(4, 5).Value = 0.1515195150002171
$ws.Cells.Item(4, 6).Value = 1.697256849473021
$ws.Cells.Item(4, 7).Value = 1.064319515806361
$ws.Cells.Item(4, 8).Value = 1.055148988853105
$ws.Cells.Item(4, 9).Value = 0.9103511746336679
$ws.Cells.Item(4, 10).Value = 0.1879787733455558
$ws.Cells.Item(4, 11).Value = 1.41170260747856
$ws.Cells.Item(4, 14).Value = 1.259080556884527

$ws.Cells.Item(5, 3).Value = 0.01730933722922856
$ws.Cells.Item(5, 4).Value = 0.1484640279873872
$ws.Cells.Item(5, 5).Value = 0.1506904007283403
$ws.Cells.Item(5, 6).Value = 1.696020123208839
$ws.Cells.Item(5, 7).Value = 1.063078267157707
$ws.Cells.Item(5, 8).Value = 1.056245968411574
$ws.Cells.Item(5, 9).Value = 0.9094762274126467
$ws.Cells.Item(5, 10).Value = 0.1869768863124861
$ws.Cells.Item(5, 11).Value = 1.374103365154212
$ws.Cells.Item(5, 14).Value = 1.2623851120021

$ws.Cells.Item(6, 3).Value = 0.01726007581734024
$ws.Cells.Item(6, 4).Value = 0.1483214323456821
$ws.Cells.Item(6, 5).Value = 0.1505539151383957
$ws.Cells.Item(6, 6).Value = 1.695832001441346
$ws.Cells.Item(6, 7).Value = 1.062885042160431
$ws.Cells.Item(6, 8).Value = 1.056437429273686
$ws.Cells.Item(6, 9).Value = 0.909341013316336
$ws.Cells.Item(6, 10).Value = 0.186812314439301
$ws.Cells.Item(6, 11).Value = 1.367865460202211
$ws.Cells.Item(6, 14).Value = 1.26294183019791

$ws.Cells.Item(7, 3).Value = 0.01760143812074944
$ws.Cells.Item(7, 4).Value = 0.1493157844040951
$ws.Cells.Item(7, 5).Value = 0.1515082536216319
$ws.Cells.Item(7, 6).Value = 1.697239014696379
$ws.Cells.Item(7, 7).Value = 1.064301911715191
$ws.Cells.Item(7, 8).Value = 1.055163158952539
$ws.Cells.Item(7, 9).Value = 0.9103386993990696
$ws.Cells.Item(7, 10).Value = 0.1879651414595429
$ws.Cells.Item(7, 11).Value = 1.411195168224481
$ws.Cells.Item(7, 14).Value = 1.259124587049556

$ws.Cells.Item(8, 3).Value = 0.01909372131397191
$ws.Cells.Item(8, 4).Value = 0.1538165937961935
$ws.Cells.Item(8, 5).Value = 0.1558922330694124
$ws.Cells.Item(8, 6).Value = 1.705957646383297
$ws.Cells.Item(8, 7).Value = 1.072432390802675
$ws.Cells.Item(8, 8).Value = 1.050894710770422
$ws.Cells.Item(8, 9).Value = 0.9162083551050131
$ws.Cells.Item(8, 10).Value = 0.1933130497587285
$ws.Cells.Item(8, 11).Value = 1.603180590583179
$ws.Cells.Item(8, 14).Value = 1.243442317106421

$ws.Cells.Item(9, 3).Value = 0.02196891395630018
$ws.Cells.Item(9, 4).Value = 0.1630449956598738
$ws.Cells.Item(9, 5).Value = 0.1651104023186818
$ws.Cells.Item(9, 6).Value = 1.732122553335415
$ws.Cells.Item(9, 7).Value = 1.09518215870726
$ws.Cells.Item(9, 8).Value = 1.04740993280241
$ws.Cells.Item(9, 9).Value = 0.9330259523735052
$ws.Cells.Item(9, 10).Value = 0.2047425792852522
$ws.Cells.Item(9, 11).Value = 1.982425150454105
$ws.Cells.Item(9, 14).Value = 1.216864277625774

$ws.Cells.Item(10, 3).Value = 0.02404985572283636
$ws.Cells.Item(10, 4).Value = 0.1700682306119603
$ws.Cells.Item(10, 5).Value = 0.1722634285003082
$ws.Cells.Item(10, 6).Value = 1.756924410543277
$ws.Cells.Item(10, 7).Value = 1.11610211174164
$ws.Cells.Item(10, 8).Value = 1.047845926592913
$ws.Cells.Item(10, 9).Value = 0.9486530990162976
$ws.Cells.Item(10, 10).Value = 0.2137206781407315
$ws.Cells.Item(10, 11).Value = 2.262738650413041
$ws.Cells.Item(10, 14).Value = 1.199880254461256

$ws.Cells.Item(11, 3).Value = 0.02498965561743205
$ws.Cells.Item(11, 4).Value = 0.1733156842200003
$ws.Cells.Item(11, 5).Value = 0.1756005775356755
$ws.Cells.Item(11, 6).Value = 1.769430050605052
$ws.Cells.Item(11, 7).Value = 1.126544858401843
$ws.Cells.Item(11, 8).Value = 1.048698887291124
$ws.Cells.Item(11, 9).Value = 0.9564805714614408
$ws.Cells.Item(11, 10).Value = 0.2179326269197901
$ws.Cells.Item(11, 11).Value = 2.390627793284182
$ws.Cells.Item(11, 14).Value = 1.192705649559265

$ws.Cells.Item(12, 3).Value = 0.02534454321223478
$ws.Cells.Item(12, 4).Value = 0.1745529167205149
$ws.Cells.Item(12, 5).Value = 0.1768762504800208
$ws.Cells.Item(12, 6).Value = 1.774342329969912
$ws.Cells.Item(12, 7).Value = 1.130633383348339
$ws.Cells.Item(12, 8).Value = 1.049116321146357
$ws.Cells.Item(12, 9).Value = 0.9595485599898126
$ws.Cells.Item(12, 10).Value = 0.2195460501556568
$ws.Cells.Item(12, 11).Value = 2.439109251901414
$ws.Cells.Item(12, 14).Value = 1.190068123222119

$ws.Cells.Item(13, 3).Value = 0.02526815622704959
$ws.Cells.Item(13, 4).Value = 0.1742861243578915
$ws.Cells.Item(13, 5).Value = 0.1766009791991507
$ws.Cells.Item(13, 6).Value = 1.773276513048643
$ws.Cells.Item(13, 7).Value = 1.129746868509301
$ws.Cells.Item(13, 8).Value = 1.049022214377743
$ws.Cells.Item(13, 9).Value = 0.9588831840669059
$ws.Cells.Item(13, 10).Value = 0.2191977488053567
$ws.Cells.Item(13, 11).Value = 2.428665578590255
$ws.Cells.Item(13, 14).Value = 1.190632632372207

$ws.Cells.Item(14, 3).Value = 0.02501887242748779
$ws.Cells.Item(14, 4).Value = 0.1734173222907032
$ws.Cells.Item(14, 5).Value = 0.1757052880612449
$ws.Cells.Item(14, 6).Value = 1.769830641682191
$ws.Cells.Item(14, 7).Value = 1.126878531385586
$ws.Cells.Item(14, 8).Value = 1.04873133544325
$ws.Cells.Item(14, 9).Value = 0.956730891300495
$ws.Cells.Item(14, 10).Value = 0.218064994017567
$ws.Cells.Item(14, 11).Value = 2.394615343451107
$ws.Cells.Item(14, 14).Value = 1.192487068347134

$ws.Cells.Item(15, 3).Value = 0.02486604907298329
$ws.Cells.Item(15, 4).Value = 0.1728861297638105
$ws.Cells.Item(15, 5).Value = 0.1751582099304727
$ws.Cells.Item(15, 6).Value = 1.767742976135438
$ws.Cells.Item(15, 7).Value = 1.125139080595289
$ws.Cells.Item(15, 8).Value = 1.048565471087755
$ws.Cells.Item(15, 9).Value = 0.9554260966635724
$ws.Cells.Item(15, 10).Value = 0.2173735539920614
$ws.Cells.Item(15, 11).Value = 2.373765412849195
$ws.Cells.Item(15, 14).Value = 1.193633297598225

$ws.Cells.Item(16, 3).Value = 0.02398829910864464
$ws.Cells.Item(16, 4).Value = 0.1698570537992907
$ws.Cells.Item(16, 5).Value = 0.1720470124103883
$ws.Cells.Item(16, 6).Value = 1.756131816770676
$ws.Cells.Item(16, 7).Value = 1.115438367409951
$ws.Cells.Item(16, 8).Value = 1.047803382060692
$ws.Cells.Item(16, 9).Value = 0.9481560632437009
$ws.Cells.Item(16, 10).Value = 0.2134479949292825
$ws.Cells.Item(16, 11).Value = 2.254388213138498
$ws.Cells.Item(16, 14).Value = 1.200360232319333

$ws.Cells.Item(17, 3).Value = 0.02344806981953695
$ws.Cells.Item(17, 4).Value = 0.1680122269671074
$ws.Cells.Item(17, 5).Value = 0.170159705162412
$ws.Cells.Item(17, 6).Value = 1.749322575814816
$ws.Cells.Item(17, 7).Value = 1.109725145814963
$ws.Cells.Item(17, 8).Value = 1.047503743988699
$ws.Cells.Item(17, 9).Value = 0.9438805825034393
$ws.Cells.Item(17, 10).Value = 0.2110725701324867
$ws.Cells.Item(17, 11).Value = 2.181248988860375
$ws.Cells.Item(17, 14).Value = 1.204628265954291

$ws.Cells.Item(18, 3).Value = 0.02313670196505768
$ws.Cells.Item(18, 4).Value = 0.1669560813047042
$ws.Cells.Item(18, 5).Value = 0.1690820070703154
$ws.Cells.Item(18, 6).Value = 1.745521162419664
$ws.Cells.Item(18, 7).Value = 1.106526171665791
$ws.Cells.Item(18, 8).Value = 1.047392998253457
$ws.Cells.Item(18, 9).Value = 0.941489043878633
$ws.Cells.Item(18, 10).Value = 0.2097183099528479
$ws.Cells.Item(18, 11).Value = 2.139216422829406
$ws.Cells.Item(18, 14).Value = 1.207135034421604

$ws.Cells.Item(19, 3).Value = 0.02303116821688178
$ws.Cells.Item(19, 4).Value = 0.1665993406523256
$ws.Cells.Item(19, 5).Value = 0.16871846194298
$ws.Cells.Item(19, 6).Value = 1.74425381099995
$ws.Cells.Item(19, 7).Value = 1.105457989440595
$ws.Cells.Item(19, 8).Value = 1.047366072042848
$ws.Cells.Item(19, 9).Value = 0.9406909041737492
$ws.Cells.Item(19, 10).Value = 0.2092618430716158
$ws.Cells.Item(19, 11).Value = 2.124990992589744
$ws.Cells.Item(19, 14).Value = 1.207992697565686

$ws.Cells.Item(20, 3).Value = 0.02350564470813055
$ws.Cells.Item(20, 4).Value = 0.1682081000988234
$ws.Cells.Item(20, 5).Value = 0.1703598017314008
$ws.Cells.Item(20, 6).Value = 1.750035515016776
$ws.Cells.Item(20, 7).Value = 1.110324304929406
$ws.Cells.Item(20, 8).Value = 1.047529263641252
$ws.Cells.Item(20, 9).Value = 0.9443287138566205
$ws.Cells.Item(20, 10).Value = 0.2113241934399923
$ws.Cells.Item(20, 11).Value = 2.189031150593792
$ws.Cells.Item(20, 14).Value = 1.204168553963164

$ws.Cells.Item(21, 3).Value = 0.02509212020033402
$ws.Cells.Item(21, 4).Value = 0.1736723076875677
$ws.Cells.Item(21, 5).Value = 0.1759680494533811
$ws.Cells.Item(21, 6).Value = 1.770837976769783
$ws.Cells.Item(21, 7).Value = 1.1277173853166
$ws.Cells.Item(21, 8).Value = 1.048814208244494
$ws.Cells.Item(21, 9).Value = 0.9573602477063901
$ws.Cells.Item(21, 10).Value = 0.218397210220914
$ws.Cells.Item(21, 11).Value = 2.404615303113985
$ws.Cells.Item(21, 14).Value = 1.19194022194862

$ws.Cells.Item(22, 3).Value = 0.02612317307313106
$ws.Cells.Item(22, 4).Value = 0.1772871247509045
$ws.Cells.Item(22, 5).Value = 0.1797031183069819
$ws.Cells.Item(22, 6).Value = 1.785463691660226
$ws.Cells.Item(22, 7).Value = 1.139866732575058
$ws.Cells.Item(22, 8).Value = 1.050204582512407
$ws.Cells.Item(22, 9).Value = 0.9664829482325104
$ws.Cells.Item(22, 10).Value = 0.2231274106424337
$ws.Cells.Item(22, 11).Value = 2.545818225290645
$ws.Cells.Item(22, 14).Value = 1.184410780905168

$ws.Cells.Item(23, 3).Value = 0.02557341541337621
$ws.Cells.Item(23, 4).Value = 0.175353856681653
$ws.Cells.Item(23, 5).Value = 0.1777032581158764
$ws.Cells.Item(23, 6).Value = 1.777563159285478
$ws.Cells.Item(23, 7).Value = 1.133310547711204
$ws.Cells.Item(23, 8).Value = 1.049412034076909
$ws.Cells.Item(23, 9).Value = 0.9615583708530551
$ws.Cells.Item(23, 10).Value = 0.2205929448364827
$ws.Cells.Item(23, 11).Value = 2.470427888019287
$ws.Cells.Item(23, 14).Value = 1.188387054857515

$ws.Cells.Item(24, 3).Value = 0.02347961753844885
$ws.Cells.Item(24, 4).Value = 0.1681195319286104
$ws.Cells.Item(24, 5).Value = 0.1702693152106249
$ws.Cells.Item(24, 6).Value = 1.749712842313755
$ws.Cells.Item(24, 7).Value = 1.110053158424989
$ws.Cells.Item(24, 8).Value = 1.047517534597318
$ws.Cells.Item(24, 9).Value = 0.9441259065979892
$ws.Cells.Item(24, 10).Value = 0.2112103990302074
$ws.Cells.Item(24, 11).Value = 2.185512784897696
$ws.Cells.Item(24, 14).Value = 1.204376224605625

$ws.Cells.Item(25, 3).Value = 0.02119660449764638
$ws.Cells.Item(25, 4).Value = 0.1605056202085677
$ws.Cells.Item(25, 5).Value = 0.1625499624409485
$ws.Cells.Item(25, 6).Value = 1.724068287865634
$ws.Cells.Item(25, 7).Value = 1.088293287914652
$ws.Cells.Item(25, 8).Value = 1.047827784326557
$ws.Cells.Item(25, 9).Value = 0.9279044817980235
$ws.Cells.Item(25, 10).Value = 0.2015490539178018
$ws.Cells.Item(25, 11).Value = 1.87953396951724
$ws.Cells.Item(25, 14).Value = 1.223607659433078
